# Update the "想去人数" (want-to-go count, column F) figures to the
# newly scraped values for the two sheets that list every event
# ("展览" = sheet 1, "全部类型" = sheet 4). The "全部类型" sheet carries
# the same rows shifted down by one (it has an extra leading row), so the
# same event lands on row N there and row N-1 on "展览".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# row on 展览 -> row on 全部类型, new F value
$updates = @(
    @{ Exhibit = 3;  All = 4;  Value = 9 },
    @{ Exhibit = 4;  All = 5;  Value = 13278 },
    @{ Exhibit = 8;  All = 9;  Value = 116 },
    @{ Exhibit = 13; All = 14; Value = 13245 },
    @{ Exhibit = 15; All = 16; Value = 579 },
    @{ Exhibit = 16; All = 17; Value = 8856 },
    @{ Exhibit = 17; All = 18; Value = 7937 },
    @{ Exhibit = 29; All = 32; Value = 124 },
    @{ Exhibit = 31; All = 34; Value = 90 }
)

foreach ($u in $updates) {
    $wsExhibit.Range("F" + $u.Exhibit).Value = $u.Value
    $wsAll.Range("F" + $u.All).Value = $u.Value
}
